$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 16, shifting the existing "Authentification" (u/login)
# and "Ajoute un utilisateur" (u/add) blocks down to rows 18 and 20.
$ws.Rows("16:17").Insert()

# Populate the newly inserted row 16 with the new "carte" list endpoint.
$ws.Range("A16:G16").Style = "20 % - Accent1"
$ws.Range("A16").Value = "c"
$ws.Range("B16").Value = "list"
$ws.Range("D16").Value = "idUser"
$ws.Range("E16").Value = "int"
$ws.Range("I16").Value = "Liste les cartes d'un utilisateur"

# Update the selection to match the saved workbook state.
$ws.Range("A17").Select()
